{"js": "// Update the date line and the table of division answers.\nconst body = context.document.body;\n\n// 1) Update the date paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.indexOf(\"2024-07-09 Tuesday\") !== -1) {\n  dateParagraph.getRange().insertText(\"2024-07-10 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the answer table contents in place (the grid size itself does\n// not change - only the cell text changes - so we overwrite the cell\n// values by their row/column position).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Expected current values, by (row, col), for the five data rows of the\n// table (rows 0, 4, 8, 12, 16 - the other rows are blank spacer rows).\nconst expectedOldByRow = {\n  0: [\"32\u00f73=10, 2\", \"34\u00f77=4, 6\", \"68\u00f77=9, 5\", \"66\u00f76=11, 0\", \"88\u00f72=44, 0\"],\n  4: [\"57\u00f72=28, 1\", \"17\u00f73=5, 2\", \"86\u00f72=43, 0\", \"15\u00f78=1, 7\", \"14\u00f78=1, 6\"],\n  8: [\"43\u00f72=21, 1\", \"51\u00f79=5, 6\", \"63\u00f78=7, 7\", \"66\u00f79=7, 3\", \"98\u00f76=16, 2\"],\n  12: [\"53\u00f73=17, 2\", \"38\u00f72=19, 0\", \"34\u00f74=8, 2\", \"97\u00f72=48, 1\", \"25\u00f72=12, 1\"],\n  16: [\"96\u00f76=16, 0\", \"88\u00f75=17, 3\", \"64\u00f77=9, 1\", \"19\u00f73=6, 1\", \"81\u00f77=11, 4\"],\n};\n\nconst newByRow = {\n  0: [\"70\u00f79=7, 7\", \"39\u00f79=4, 3\", \"99\u00f75=19, 4\", \"22\u00f75=4, 2\", \"23\u00f72=11, 1\"],\n  4: [\"91\u00f74=22, 3\", \"55\u00f75=11, 0\", \"79\u00f79=8, 7\", \"32\u00f76=5, 2\", \"62\u00f77=8, 6\"],\n  8: [\"63\u00f79=7, 0\", \"98\u00f76=16, 2\", \"23\u00f79=2, 5\", \"58\u00f79=6, 4\", \"21\u00f79=2, 3\"],\n  12: [\"56\u00f73=18, 2\", \"17\u00f77=2, 3\", \"13\u00f75=2, 3\", \"47\u00f73=15, 2\", \"67\u00f74=16, 3\"],\n  16: [\"30\u00f73=10, 0\", \"36\u00f76=6, 0\", \"53\u00f76=8, 5\", \"64\u00f77=9, 1\", \"28\u00f79=3, 1\"],\n};\n\nconst values = table.values;\nfor (const rowIndexStr of Object.keys(newByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  if (rowIndex >= values.length) continue;\n  const expectedOld = expectedOldByRow[rowIndex];\n  const newRow = newByRow[rowIndex];\n  for (let c = 0; c < newRow.length; c++) {\n    if (c >= values[rowIndex].length) continue;\n    // Only overwrite when the existing text matches what we expect, so we\n    // never clobber unrelated content if the table shape differs.\n    if (values[rowIndex][c] === expectedOld[c]) {\n      values[rowIndex][c] = newRow[c];\n    }\n  }\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Update the date line and the table of division answers.\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the body).\n$p1 = $d.Paragraphs.Item(1)\n$p1Text = $p1.Range.Text.TrimEnd([char]13)\nif ($p1Text -eq \"2024-07-09 Tuesday\") {\n    $p1.Range.Text = \"2024-07-10 Wednesday\"\n}\n\n# 2) Update the answer table contents in place (the grid size itself does\n# not change - only the cell text changes).\n$t = $d.Tables.Item(1)\n\n# Each entry is: row, column, expected old text, new text.\n$updates = @(\n    @(1, 1, \"32\u00f73=10, 2\", \"70\u00f79=7, 7\"),\n    @(1, 2, \"34\u00f77=4, 6\", \"39\u00f79=4, 3\"),\n    @(1, 3, \"68\u00f77=9, 5\", \"99\u00f75=19, 4\"),\n    @(1, 4, \"66\u00f76=11, 0\", \"22\u00f75=4, 2\"),\n    @(1, 5, \"88\u00f72=44, 0\", \"23\u00f72=11, 1\"),\n\n    @(5, 1, \"57\u00f72=28, 1\", \"91\u00f74=22, 3\"),\n    @(5, 2, \"17\u00f73=5, 2\", \"55\u00f75=11, 0\"),\n    @(5, 3, \"86\u00f72=43, 0\", \"79\u00f79=8, 7\"),\n    @(5, 4, \"15\u00f78=1, 7\", \"32\u00f76=5, 2\"),\n    @(5, 5, \"14\u00f78=1, 6\", \"62\u00f77=8, 6\"),\n\n    @(9, 1, \"43\u00f72=21, 1\", \"63\u00f79=7, 0\"),\n    @(9, 2, \"51\u00f79=5, 6\", \"98\u00f76=16, 2\"),\n    @(9, 3, \"63\u00f78=7, 7\", \"23\u00f79=2, 5\"),\n    @(9, 4, \"66\u00f79=7, 3\", \"58\u00f79=6, 4\"),\n    @(9, 5, \"98\u00f76=16, 2\", \"21\u00f79=2, 3\"),\n\n    @(13, 1, \"53\u00f73=17, 2\", \"56\u00f73=18, 2\"),\n    @(13, 2, \"38\u00f72=19, 0\", \"17\u00f77=2, 3\"),\n    @(13, 3, \"34\u00f74=8, 2\", \"13\u00f75=2, 3\"),\n    @(13, 4, \"97\u00f72=48, 1\", \"47\u00f73=15, 2\"),\n    @(13, 5, \"25\u00f72=12, 1\", \"67\u00f74=16, 3\"),\n\n    @(17, 1, \"96\u00f76=16, 0\", \"30\u00f73=10, 0\"),\n    @(17, 2, \"88\u00f75=17, 3\", \"36\u00f76=6, 0\"),\n    @(17, 3, \"64\u00f77=9, 1\", \"53\u00f76=8, 5\"),\n    @(17, 4, \"19\u00f73=6, 1\", \"64\u00f77=9, 1\"),\n    @(17, 5, \"81\u00f77=11, 4\", \"28\u00f79=3, 1\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $oldVal = $u[2]\n    $newVal = $u[3]\n\n    $cell = $t.Cell($row, $col)\n    $cellText = $cell.Range.Text.TrimEnd([char]7).TrimEnd([char]13)\n\n    if ($cellText -eq $oldVal) {\n        $cell.Range.Text = $newVal\n    }\n}\n"}
